$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns stay text (avoid Excel auto-converting
# numeric-looking strings like "47.00" or "0.999" into Number cells, which
# would drop formatting such as trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '40.098.76'
$ws.Range("E2").Value = '  +3.67%  '
$ws.Range("D3").Value = '2.238.58'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '294.87'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = '86.86'
$ws.Range("E6").Value = '  +8.45%  '
$ws.Range("E7").Value = '  +2.95%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +4.18%  '
$ws.Range("D10").Value = '31.34'
$ws.Range("E10").Value = '  +12.98%  '
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").Value = '47.00'
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("E14").Value = '  +6.22%  '
$ws.Range("D15").Value = '2.578.20'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '14.17'
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("D17").Value = '2.268.42'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("D18").Value = '0.733'
$ws.Range("E18").Value = '  +4.50%  '
$ws.Range("D19").Value = '40.027.03'
$ws.Range("E19").Value = '  +3.73%  '
$ws.Range("D20").Value = '0.0₃0892'
$ws.Range("E20").Value = '  +4.20%  '
$ws.Range("E21").Value = '  +3.59%  '
$ws.Range("D22").Value = '10.92'
$ws.Range("E22").Value = '  +11.63%  '
$ws.Range("D23").Value = '65.25'
$ws.Range("E23").Value = '  +2.37%  '
$ws.Range("D24").Value = '235.63'
$ws.Range("E24").Value = '  +6.28%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  +4.92%  '
$ws.Range("E27").Value = '  +7.37%  '
$ws.Range("E28").Value = '  +3.66%  '
$ws.Range("D30").Value = '9.26'
$ws.Range("E30").Value = '  +4.56%  '
$ws.Range("D31").Value = '33.42'
$ws.Range("E31").Value = '  +8.07%  '
$ws.Range("D32").Value = '152.47'
$ws.Range("E32").Value = '  +4.07%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '4.89'
$ws.Range("E34").Value = '  +3.70%  '
$ws.Range("E35").Value = '  +5.21%  '
$ws.Range("E36").Value = '  +3.67%  '
$ws.Range("D37").Value = '16.33'
$ws.Range("E37").Value = '  +15.52%  '
$ws.Range("E38").Value = '  +3.33%  '
$ws.Range("E39").Value = '  +6.79%  '
$ws.Range("D40").Value = '2.75'
$ws.Range("E40").Value = '  +6.20%  '
$ws.Range("D41").Value = '1.71'
$ws.Range("E41").Value = '  +8.39%  '
$ws.Range("E42").Value = '  +6.79%  '
$ws.Range("D43").Value = '2.045.79'
$ws.Range("E43").Value = '  +9.02%  '
$ws.Range("E44").Value = '  +7.89%  '
$ws.Range("D45").Value = '10.10'
$ws.Range("E45").Value = '  +14.21%  '
$ws.Range("E46").Value = '  +6.84%  '
$ws.Range("D47").Value = '16.49'
$ws.Range("E47").Value = '  +5.43%  '
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").Value = '2.448.90'
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("D50").Value = '70.82'
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '89.59'
$ws.Range("E51").Value = '  +5.38%  '
